# Update cryptos list prices and Volume(1h) percentages.
# Mirrors the upstream GitHub Actions scrape-refresh commit.
#
# Note: several Price (column D) values look like plain decimals
# (e.g. "585.67") which Excel would normally auto-convert to a
# number on assignment, losing the exact text formatting (trailing
# zeros, exact digits). We force those through as literal text via
# the leading apostrophe text-qualifier, the same mechanism a user
# typing into Excel would use, then reset the cell Style back to
# "Normal" so no extra number-format styling is left behind on the
# cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.786.61"
$ws.Range("E2").Value = "  -3.45%  "
$ws.Range("D3").Value = "2.901.36"
$ws.Range("E3").Value = "  -4.73%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'585.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "'145.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.83%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.501"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "2.898.19"
$ws.Range("E9").Value = "  -4.77%  "
$ws.Range("D10").Value = "'6.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.88%  "
$ws.Range("E11").Value = "  -6.16%  "
$ws.Range("D12").Value = "'0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.65%  "
$ws.Range("E13").Value = "  -5.26%  "
$ws.Range("D14").Value = "'33.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.81%  "
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "3.379.56"
$ws.Range("E16").Value = "  -4.73%  "
$ws.Range("D17").Value = "60.674.59"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").Value = "'6.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.54%  "
$ws.Range("D19").Value = "2.900.44"
$ws.Range("E19").Value = "  -4.51%  "
$ws.Range("D20").Value = "'423.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.58%  "
$ws.Range("D21").Value = "'13.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.31%  "
$ws.Range("D22").Value = "'0.669"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("D23").Value = "'7.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.01%  "
$ws.Range("D24").Value = "'79.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.31%  "
$ws.Range("D25").Value = "'10.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  -3.71%  "
$ws.Range("D27").Value = "'11.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'7.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("E31").Value = "  -4.26%  "
$ws.Range("D32").Value = "'2.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "'26.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.20%  "
$ws.Range("D34").Value = "'0.106"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.90%  "
$ws.Range("D35").Value = "0.0₃0834"
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("D37").Value = "'5.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.25%  "
$ws.Range("D38").Value = "'49.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.31%  "
$ws.Range("E39").Value = "  -6.78%  "
$ws.Range("E40").Value = "  -4.90%  "
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").Value = "'8.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.10%  "
$ws.Range("D43").Value = "'0.290"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").Value = "'41.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "'0.0346"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("D46").Value = "'372.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.70%  "
$ws.Range("D47").Value = "2.663.41"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").Value = "'133.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "'24.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.37%  "
$ws.Range("E51").Value = "  -2.33%  "
